# Updated cryptos list on Thu Jan  4 13:43:41 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: price values in column D are plain text (not real numbers) in the
# source workbook (e.g. "317.70", "7.69", thousand-grouped "43.328.52", the
# shiba-inu subscript notation, etc.). A leading apostrophe is used so the
# COM layer stores them as text instead of re-parsing/rounding them as
# numeric values.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'43.353.10"
$ws.Range("E2").Value = "  +1.22%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.235.83"
$ws.Range("E3").Value = "  +0.18%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'317.50"
$ws.Range("E5").Value = "  +0.68%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'99.53"
$ws.Range("E6").Value = "  -0.29%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +1.80%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.01%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.09%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'37.08"
$ws.Range("E10").Value = "  -0.68%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "'0.0831"
$ws.Range("E11").Value = "  -0.84%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "'7.70"
$ws.Range("E12").Value = "  +0.33%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +1.78%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "'0.864"
$ws.Range("E14").Value = "  -1.39%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'14.35"
$ws.Range("E15").Value = "  +2.73%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "'2.244.54"
$ws.Range("E16").Value = "  +1.27%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "'43.288.64"
$ws.Range("E17").Value = "  +1.52%  "

# Row 18 - InternetComputer(DFINITY)
$ws.Range("D18").Value = "'14.16"
$ws.Range("E18").Value = "  -0.15%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "'6.66"
$ws.Range("E19").Value = "  +0.09%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "'0.0" + [char]0x2083 + "0973"
$ws.Range("E20").Value = "  +2.42%  "

# Row 21 - Litecoin
$ws.Range("D21").Value = "'65.33"
$ws.Range("E21").Value = "  +0.39%  "

# Row 22 - PancakeSwap
$ws.Range("E22").Value = "  -3.03%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'236.19"
$ws.Range("E23").Value = "  +0.11%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  +1.90%  "

# Row 25 - Dai
$ws.Range("E25").Value = "  -0.20%  "

# Row 26 - LEO
$ws.Range("D26").Value = "'4.06"
$ws.Range("E26").Value = "  +2.96%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "'10.05"
$ws.Range("E27").Value = "  -1.90%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  +1.98%  "

# Row 29 / Row 30 - InjectiveProtocol and Filecoin swap places
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'6.38"
$ws.Range("E29").Value = "  -3.51%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'36.63"
$ws.Range("E30").Value = "  +9.71%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "'20.27"
$ws.Range("E31").Value = "  -1.44%  "

# Row 32 - Hedera
$ws.Range("D32").Value = "'0.0872"
$ws.Range("E32").Value = "  -2.27%  "

# Row 33 - Monero
$ws.Range("D33").Value = "'157.40"
$ws.Range("E33").Value = "  -1.59%  "

# Row 34 - WEMIXToken
$ws.Range("E34").Value = "  -0.90%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "'3.18"
$ws.Range("E35").Value = "  +1.59%  "

# Row 36 - Stellar
$ws.Range("E36").Value = "  -1.65%  "

# Row 37 - ARBITRUM
$ws.Range("E37").Value = "  +0.53%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  -2.19%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +1.02%  "

# Row 40 - NEARProtocol
$ws.Range("D40").Value = "'3.70"
$ws.Range("E40").Value = "  +2.84%  "

# Row 41 - VeChain
$ws.Range("E41").Value = "  -0.86%  "

# Row 42 - Celestia
$ws.Range("D42").Value = "'14.35"
$ws.Range("E42").Value = "  +18.17%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.07%  "

# Row 44 - Maker
$ws.Range("D44").Value = "'1.822.44"
$ws.Range("E44").Value = "  -0.30%  "

# Row 45 - Algorand
$ws.Range("D45").Value = "'0.203"
$ws.Range("E45").Value = "  -2.39%  "

# Row 46 - BitcoinSV
$ws.Range("D46").Value = "'83.94"
$ws.Range("E46").Value = "  -6.75%  "

# Row 47 - THORChain
$ws.Range("E47").Value = "  -2.23%  "

# Row 48 - FraxShare
$ws.Range("D48").Value = "'8.82"
$ws.Range("E48").Value = "  +2.44%  "

# Row 49 - ordi
$ws.Range("D49").Value = "'73.84"
$ws.Range("E49").Value = "  -5.65%  "

# Row 50 - Aave
$ws.Range("D50").Value = "'103.21"
$ws.Range("E50").Value = "  +0.77%  "

# Row 51 - MultiversX
$ws.Range("D51").Value = "'58.13"
$ws.Range("E51").Value = "  -4.94%  "
